$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Cells.Item(2, 2).Value = 77388.66797673714
$ws.Cells.Item(2, 4).Value = 9992.97670278544
$ws.Cells.Item(2, 5).Value = 1770
$ws.Cells.Item(2, 6).Value = 38945.6536412876

$ws = $wb.Worksheets.Item("Capacities")
$ws.Cells.Item(3, 3).Value = 68

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Cells.Item(2, 7).Value = 13.6
$ws.Cells.Item(2, 8).Value = 27.2
$ws.Cells.Item(2, 9).Value = 34
$ws.Cells.Item(2, 10).Value = 40.8
$ws.Cells.Item(2, 11).Value = 47.6
$ws.Cells.Item(2, 12).Value = 54.4
$ws.Cells.Item(2, 13).Value = 61.2
$ws.Cells.Item(2, 14).Value = 68
$ws.Cells.Item(2, 15).Value = 61.2
$ws.Cells.Item(2, 16).Value = 54.4
$ws.Cells.Item(2, 17).Value = 47.6
$ws.Cells.Item(2, 18).Value = 34
$ws.Cells.Item(2, 19).Value = 20.4
$ws.Cells.Item(2, 20).Value = 13.6
$ws.Cells.Item(3, 9).Value = 27.2
$ws.Cells.Item(3, 10).Value = 40.8
$ws.Cells.Item(3, 11).Value = 54.4
$ws.Cells.Item(3, 12).Value = 61.2
$ws.Cells.Item(3, 13).Value = 68
$ws.Cells.Item(3, 14).Value = 54.4
$ws.Cells.Item(3, 15).Value = 47.6
$ws.Cells.Item(3, 16).Value = 34
$ws.Cells.Item(3, 17).Value = 34
$ws.Cells.Item(3, 18).Value = 20.4
$ws.Cells.Item(3, 19).Value = 13.6
$ws.Cells.Item(4, 11).Value = 27.2
$ws.Cells.Item(4, 12).Value = 47.6
$ws.Cells.Item(4, 13).Value = 54.4
$ws.Cells.Item(4, 14).Value = 54.4
$ws.Cells.Item(4, 15).Value = 47.6
$ws.Cells.Item(4, 16).Value = 27.2
$ws.Cells.Item(4, 17).Value = 10.38312417100186

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Cells.Item(2, 7).Value = 64.3
$ws.Cells.Item(2, 8).Value = 14.2
$ws.Cells.Item(2, 9).Value = 2.8
$ws.Cells.Item(2, 10).Value = 1.8
$ws.Cells.Item(2, 11).Value = 21.6
$ws.Cells.Item(2, 12).Value = 33.6
$ws.Cells.Item(2, 13).Value = 37.8
$ws.Cells.Item(2, 14).Value = 42
$ws.Cells.Item(2, 15).Value = 30
$ws.Cells.Item(2, 16).Value = 25.8
$ws.Cells.Item(2, 17).Value = 132.0342720130611
$ws.Cells.Item(2, 18).Value = 0.2
$ws.Cells.Item(2, 19).Value = 32.4
$ws.Cells.Item(2, 20).Value = 45.6
$ws.Cells.Item(3, 9).Value = 27.43079277624771
$ws.Cells.Item(3, 10).Value = 40.8
$ws.Cells.Item(3, 11).Value = 54.4
$ws.Cells.Item(3, 12).Value = 61.2
$ws.Cells.Item(3, 13).Value = 44.6
$ws.Cells.Item(3, 14).Value = 28.4
$ws.Cells.Item(3, 15).Value = 47.6
$ws.Cells.Item(3, 16).Value = 5.4
$ws.Cells.Item(3, 17).Value = 8
$ws.Cells.Item(3, 18).Value = 20.4
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 27.2
$ws.Cells.Item(4, 12).Value = 47.6
$ws.Cells.Item(4, 13).Value = 31
$ws.Cells.Item(4, 14).Value = 54.4
$ws.Cells.Item(4, 15).Value = 47.6
$ws.Cells.Item(4, 16).Value = 27.2
$ws.Cells.Item(4, 17).Value = 10.38312417100186
$ws.Cells.Item(4, 18).Value = 0

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Cells.Item(2, 7).Value = 183.657
$ws.Cells.Item(2, 8).Value = 197.715
$ws.Cells.Item(2, 9).Value = 200.487
$ws.Cells.Item(2, 10).Value = 202.269
$ws.Cells.Item(2, 11).Value = 223.653
$ws.Cells.Item(2, 12).Value = 256.917
$ws.Cells.Item(2, 13).Value = 294.3389999999999
$ws.Cells.Item(2, 14).Value = 335.9189999999999
$ws.Cells.Item(2, 15).Value = 365.6189999999999
$ws.Cells.Item(2, 16).Value = 391.1609999999999
$ws.Cells.Item(2, 17).Value = 521.8749292929305
$ws.Cells.Item(2, 18).Value = 522.0729292929304
$ws.Cells.Item(2, 19).Value = 554.1489292929305
$ws.Cells.Item(3, 9).Value = 147.1564848484852
$ws.Cells.Item(3, 10).Value = 187.5484848484852
$ws.Cells.Item(3, 11).Value = 241.4044848484852
$ws.Cells.Item(3, 12).Value = 301.9924848484852
$ws.Cells.Item(3, 13).Value = 346.1464848484852
$ws.Cells.Item(3, 14).Value = 374.2624848484852
$ws.Cells.Item(3, 15).Value = 421.3864848484852
$ws.Cells.Item(3, 16).Value = 426.7324848484852
$ws.Cells.Item(3, 17).Value = 434.6524848484852
$ws.Cells.Item(3, 18).Value = 454.8484848484852
$ws.Cells.Item(4, 10).Value = 120
$ws.Cells.Item(4, 11).Value = 146.928
$ws.Cells.Item(4, 12).Value = 194.052
$ws.Cells.Item(4, 13).Value = 224.742
$ws.Cells.Item(4, 14).Value = 278.598
$ws.Cells.Item(4, 15).Value = 325.722
$ws.Cells.Item(4, 16).Value = 352.65
$ws.Cells.Item(4, 17).Value = 362.9292929292918

$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 17).Value = 110.4342720130611
$ws.Cells.Item(2, 20).Value = 34
$ws.Cells.Item(3, 9).Value = 0.2307927762477106
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(4, 16).Value = 0

$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0

$ws = $wb.Worksheets.Item("Feed in from Type 4")
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 19).Value = 9.6
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0

Write-Output "Applied all Year 5 test-case edits"